$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell-value corrections as described by the target diff.
# Each group below corresponds to the <c> cell value changes
# in xl/worksheets/sheet1.xml for a single row (stock/value
# recalculations and a couple of adjacent data rows whose
# contents were swapped).

$ws.Range("F132").Value = 22
$ws.Range("G132").Value = 1088.56

$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0

$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0

$ws.Range("B138").Value = 1056.54

$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0

$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0

$ws.Range("B184").Value = 3465.8

$ws.Range("F235").Value = 55
$ws.Range("G235").Value = 4435.2

$ws.Range("F237").Value = 38
$ws.Range("G237").Value = 1898.86

$ws.Range("F238").Value = 56
$ws.Range("G238").Value = 3358.32

$ws.Range("F244").Value = 15
$ws.Range("G244").Value = 499.5

$ws.Range("B246").Value = 48706
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2

$ws.Range("B247").Value = 64973
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 64
$ws.Range("G247").Value = 2131.2

$ws.Range("F249").Value = 13
$ws.Range("G249").Value = 632.97

$ws.Range("F250").Value = 32
$ws.Range("G250").Value = 3665.92

$ws.Range("F263").Value = 6
$ws.Range("G263").Value = 582.6

$ws.Range("B274").Value = 89778.35000000001

$ws.Range("B277").Value = 63565
$ws.Range("E277").Value = 109.19
$ws.Range("F277").Value = 60
$ws.Range("G277").Value = 6162.6

$ws.Range("B278").Value = 61610
$ws.Range("E278").Value = 122.71
$ws.Range("F278").Value = -58
$ws.Range("G278").Value = -5957.18

$ws.Range("B294").Value = 63531
$ws.Range("F294").Value = 80
$ws.Range("G294").Value = 11478.4

$ws.Range("B296").Value = 63571
$ws.Range("F296").Value = 4
$ws.Range("G296").Value = 573.92

$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 145
$ws.Range("G299").Value = 6907.8

$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12

$ws.Range("B311").Value = 63563
$ws.Range("E311").Value = 119.04
$ws.Range("F311").Value = 2
$ws.Range("G311").Value = 223.92

$ws.Range("B312").Value = 61605
$ws.Range("E312").Value = 133.78
$ws.Range("F312").Value = -13
$ws.Range("G312").Value = -1455.48

$ws.Range("B356").Value = 63681
$ws.Range("E356").Value = 23.84
$ws.Range("F356").Value = 0
$ws.Range("G356").Value = 0

$ws.Range("B357").Value = 31930
$ws.Range("E357").Value = 26.8
$ws.Range("F357").Value = -62
$ws.Range("G357").Value = -1390.04

$ws.Range("F405").Value = 26
$ws.Range("G405").Value = 210.34

$ws.Range("F406").Value = 0
$ws.Range("G406").Value = 0

$ws.Range("B411").Value = 3177.98

$ws.Range("F413").Value = 0
$ws.Range("G413").Value = 0

$ws.Range("B416").Value = 919.98

$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 42
$ws.Range("G420").Value = 4432.68

$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2

$ws.Range("B467").Value = 65068
$ws.Range("E467").Value = 13.97
$ws.Range("F467").Value = 113
$ws.Range("G467").Value = 1485.95

$ws.Range("B468").Value = 53602
$ws.Range("E468").Value = 15.69
$ws.Range("F468").Value = -231
$ws.Range("G468").Value = -3037.65

$ws.Range("B476").Value = 45706
$ws.Range("E476").Value = 23.58
$ws.Range("F476").Value = -202
$ws.Range("G476").Value = -3985.46

$ws.Range("B477").Value = 64922
$ws.Range("E477").Value = 20.98
$ws.Range("F477").Value = 117
$ws.Range("G477").Value = 2308.41

$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 203
$ws.Range("G485").Value = 2669.45

$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945

$ws.Range("F566").Value = 0
$ws.Range("G566").Value = 0

$ws.Range("B568").Value = 138.54

$ws.Range("B603").Value = 60031
$ws.Range("E603").Value = 111.69
$ws.Range("F603").Value = -5
$ws.Range("G603").Value = -492.5

$ws.Range("B604").Value = 64836
$ws.Range("E604").Value = 104.71
$ws.Range("F604").Value = 3
$ws.Range("G604").Value = 295.5

$ws.Range("B717").Value = 63150
$ws.Range("D717").Value = 75.68000000000001
$ws.Range("E717").Value = 80.45
$ws.Range("F717").Value = 68
$ws.Range("G717").Value = 5146.24

$ws.Range("B718").Value = 61428
$ws.Range("D718").Value = 69.16
$ws.Range("E718").Value = 73.52
$ws.Range("F718").Value = 1
$ws.Range("G718").Value = 69.16

$ws.Range("F730").Value = 18
$ws.Range("G730").Value = 4490.64

$ws.Range("B744").Value = 65079
$ws.Range("F744").Value = 21
$ws.Range("G744").Value = 858.27

$ws.Range("B745").Value = 65362
$ws.Range("F745").Value = 54
$ws.Range("G745").Value = 2206.98

$ws.Range("B755").Value = 84454.02

$ws.Range("F781").Value = 63
$ws.Range("G781").Value = 11098.71

$ws.Range("B787").Value = 796288.33

$ws.Range("B805").Value = 3106076.03

$ws.Range("B806").Value = 3106076.03

